# Bug fix on the formulas for teams
#
# 1) The pool draw on the "data" sheet is re-shuffled (player/team pairs
#    moved between pools while the Pool-header cells in column A keep
#    their original position/order).
# 2) Several "Elimination Matches" formulas had a broken/blank sheet
#    reference (CONCATENATE("Pool X.Y ",''!)) - these are restored to
#    point at the correct cell on the "Pool Matches" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Re-draw the pools on the "data" sheet (columns B = Player Name,
#    C = Player Dojo/Team). Column A (pool headers) is unchanged.
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("data")

$ws.Range("B3").Value = "Legolas Greenleaf"
$ws.Range("C3").Value = "Team Mu"

$ws.Range("B5").Value = "Voldemort"
$ws.Range("C5").Value = "Team Chi"

$ws.Range("B6").Value = "Othello"
$ws.Range("C6").Value = "Team Omicron"

$ws.Range("B7").Value = "Hermione Granger"
$ws.Range("C7").Value = "Team Theta"

$ws.Range("B8").Value = "Tyrion Lannister"
$ws.Range("C8").Value = "Team Upsilon"

$ws.Range("B9").Value = "Frodo Baggins"
$ws.Range("C9").Value = "Team Zeta"

$ws.Range("B10").Value = "Neville Longbottom"
$ws.Range("C10").Value = "Team Xi"

$ws.Range("B11").Value = "Ygritte"
$ws.Range("C11").Value = "Team Alpha"

$ws.Range("B12").Value = "Petyr Baelish"
$ws.Range("C12").Value = "Team Pi"

$ws.Range("B13").Value = "Quirinus Quirrell"
$ws.Range("C13").Value = "Team Rho"

$ws.Range("B14").Value = "Katniss Everdeen"
$ws.Range("C14").Value = "Team Lambda"

$ws.Range("B15").Value = "Jon Snow"
$ws.Range("C15").Value = "Team Kappa"

$ws.Range("B16").Value = "Daenerys Targaryen"
$ws.Range("C16").Value = "Team Delta"

$ws.Range("B17").Value = "Xaro Xhoan Daxos"
$ws.Range("C17").Value = "Team Omega"

$ws.Range("B18").Value = "Willy Wonka"
$ws.Range("C18").Value = "Team Psi"

$ws.Range("B19").Value = "Moby Dick"
$ws.Range("C19").Value = "Team Nu"

$ws.Range("B20").Value = "Ron Weasley"
$ws.Range("C20").Value = "Team Sigma"

$ws.Range("B21").Value = "Cersei Lannister"
$ws.Range("C21").Value = "Team Gamma"

$ws.Range("B22").Value = "Eddard Stark"
$ws.Range("C22").Value = "Team Epsilon"

$ws.Range("B23").Value = "Gandalf The Grey"
$ws.Range("C23").Value = "Team Eta"

$ws.Range("B24").Value = "Inigo Montoya"
$ws.Range("C24").Value = "Team Iota"

# ---------------------------------------------------------------------
# 2) Fix the broken CONCATENATE(...,''!) formulas on "Elimination
#    Matches" so they point at the right "Pool Matches" cell.
# ---------------------------------------------------------------------
$em = $wb.Worksheets.Item("Elimination Matches")

$em.Range("G5").Formula  = "=CONCATENATE(""Pool G.2 "",'Pool Matches'!G171)"
$em.Range("G10").Formula = "=CONCATENATE(""Pool G.2 "",'Pool Matches'!G171)"

$em.Range("O5").Formula  = "=CONCATENATE(""Pool F.2 "",'Pool Matches'!O123)"
$em.Range("O10").Formula = "=CONCATENATE(""Pool F.2 "",'Pool Matches'!O123)"

$em.Range("G19").Formula = "=CONCATENATE(""Pool E.2 "",'Pool Matches'!G123)"
$em.Range("G24").Formula = "=CONCATENATE(""Pool E.2 "",'Pool Matches'!G123)"

$em.Range("I19").Formula = "=CONCATENATE(""Pool B.2 "",'Pool Matches'!O37)"
$em.Range("I24").Formula = "=CONCATENATE(""Pool B.2 "",'Pool Matches'!O37)"

$em.Range("O19").Formula = "=CONCATENATE(""Pool A.2 "",'Pool Matches'!G48)"
$em.Range("O24").Formula = "=CONCATENATE(""Pool A.2 "",'Pool Matches'!G48)"

$em.Range("G33").Formula = "=CONCATENATE(""Pool C.2 "",'Pool Matches'!G86)"
$em.Range("G38").Formula = "=CONCATENATE(""Pool C.2 "",'Pool Matches'!G86)"

$em.Range("O33").Formula = "=CONCATENATE(""Pool D.2 "",'Pool Matches'!O86)"
$em.Range("O38").Formula = "=CONCATENATE(""Pool D.2 "",'Pool Matches'!O86)"
